$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2157434402332362
$ws.Range("C2").Value = 0.5043731778425656
$ws.Range("J2").Value = 0.01166180758017493
$ws.Range("P2").Value = 0.163265306122449
$ws.Range("S2").Value = 0.1049562682215743
$ws.Range("B3").Value = 0.005555555555555556
$ws.Range("C3").Value = 0.03888888888888889
$ws.Range("J3").Value = 0.02777777777777778
$ws.Range("P3").Value = 0.6777777777777778
$ws.Range("S3").Value = 0.25
$ws.Range("J4").Value = 0.1111111111111111
$ws.Range("P4").Value = 0.6666666666666666
$ws.Range("S4").Value = 0.2222222222222222
$ws.Range("B6").Value = 0.05905511811023622
$ws.Range("D6").Value = 0.02362204724409449
$ws.Range("F6").Value = 0.04724409448818898
$ws.Range("J6").Value = 0.2559055118110236
$ws.Range("O6").Value = 0.01181102362204724
$ws.Range("Q6").Value = 0.1535433070866142
$ws.Range("R6").Value = 0.07874015748031496
$ws.Range("S6").Value = 0.3700787401574803
$ws.Range("B7").Value = 0.1275510204081633
$ws.Range("D7").Value = 0.01020408163265306
$ws.Range("F7").Value = 0.06122448979591837
$ws.Range("J7").Value = 0.1377551020408163
$ws.Range("O7").Value = 0.00510204081632653
$ws.Range("Q7").Value = 0.1326530612244898
$ws.Range("R7").Value = 0.08163265306122448
$ws.Range("S7").Value = 0.4438775510204082
$ws.Range("B8").Value = 0.1038374717832957
$ws.Range("D8").Value = 0.02031602708803612
$ws.Range("E8").Value = 0.002257336343115124
$ws.Range("F8").Value = 0.06320541760722348
$ws.Range("J8").Value = 0.09932279909706546
$ws.Range("O8").Value = 0.01128668171557562
$ws.Range("Q8").Value = 0.1738148984198646
$ws.Range("R8").Value = 0.09932279909706546
$ws.Range("S8").Value = 0.4266365688487584
$ws.Range("B9").Value = 0.07239819004524888
$ws.Range("D9").Value = 0.03167420814479638
$ws.Range("F9").Value = 0.08597285067873303
$ws.Range("J9").Value = 0.05882352941176471
$ws.Range("O9").Value = 0.03167420814479638
$ws.Range("Q9").Value = 0.253393665158371
$ws.Range("R9").Value = 0.07692307692307693
$ws.Range("S9").Value = 0.3891402714932127
$ws.Range("B10").Value = 0.1185682326621924
$ws.Range("D10").Value = 0.01715137956748695
$ws.Range("E10").Value = 0.0007457121551081282
$ws.Range("F10").Value = 0.07307979120059657
$ws.Range("J10").Value = 0.0842654735272185
$ws.Range("O10").Value = 0.01938851603281133
$ws.Range("Q10").Value = 0.1730052199850858
$ws.Range("R10").Value = 0.09619686800894854
$ws.Range("S10").Value = 0.4175988068605518
$ws.Range("F11").Value = 0.003164556962025316
$ws.Range("G11").Value = 0.1329113924050633
$ws.Range("J11").Value = 0.1139240506329114
$ws.Range("K11").Value = 0.189873417721519
$ws.Range("L11").Value = 0.5443037974683544
$ws.Range("S11").Value = 0.01582278481012658
$ws.Range("G12").Value = 0.7965116279069767
$ws.Range("J12").Value = 0.1511627906976744
$ws.Range("K12").Value = 0.01162790697674419
$ws.Range("L12").Value = 0.005813953488372093
$ws.Range("S12").Value = 0.03488372093023256
$ws.Range("G13").Value = 0.5625
$ws.Range("J13").Value = 0.3125
$ws.Range("S13").Value = 0.125
$ws.Range("F14").Value = 0.25
$ws.Range("G14").Value = 0.75
$ws.Range("F15").Value = 0.05
$ws.Range("H15").Value = 0.15
$ws.Range("I15").Value = 0.1125
$ws.Range("J15").Value = 0.35
$ws.Range("K15").Value = 0.05416666666666667
$ws.Range("M15").Value = 0.0125
$ws.Range("O15").Value = 0.03333333333333333
$ws.Range("S15").Value = 0.2375
$ws.Range("F16").Value = 0.025
$ws.Range("H16").Value = 0.175
$ws.Range("I16").Value = 0.095
$ws.Range("J16").Value = 0.385
$ws.Range("K16").Value = 0.12
$ws.Range("M16").Value = 0.02
$ws.Range("O16").Value = 0.065
$ws.Range("S16").Value = 0.115
$ws.Range("F17").Value = 0.004705882352941176
$ws.Range("H17").Value = 0.1694117647058823
$ws.Range("I17").Value = 0.09411764705882353
$ws.Range("J17").Value = 0.4658823529411765
$ws.Range("K17").Value = 0.08470588235294117
$ws.Range("M17").Value = 0.01411764705882353
$ws.Range("N17").Value = 0.002352941176470588
$ws.Range("O17").Value = 0.05176470588235294
$ws.Range("S17").Value = 0.1129411764705882
$ws.Range("F18").Value = 0.02232142857142857
$ws.Range("H18").Value = 0.1607142857142857
$ws.Range("I18").Value = 0.07142857142857142
$ws.Range("J18").Value = 0.4642857142857143
$ws.Range("K18").Value = 0.08035714285714286
$ws.Range("M18").Value = 0.004464285714285714
$ws.Range("O18").Value = 0.07589285714285714
$ws.Range("S18").Value = 0.1205357142857143
$ws.Range("F19").Value = 0.0226628895184136
$ws.Range("H19").Value = 0.1890934844192635
$ws.Range("I19").Value = 0.08640226628895184
$ws.Range("J19").Value = 0.3930594900849859
$ws.Range("K19").Value = 0.108356940509915
$ws.Range("M19").Value = 0.01345609065155807
$ws.Range("N19").Value = 0.002124645892351275
$ws.Range("O19").Value = 0.07719546742209632
$ws.Range("S19").Value = 0.1076487252124646